$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix data errors ("sua loi du lieu") ---

# Account (A2) password value was wrong -> fix it
$ws.Range("A2").Value = "kikaho"

# The old row 3 (A3/B3) held a stray second account; its data is wrong.
# Replace it with three corrected values placed alongside row 2 instead
# (C2:E2), then remove the now-redundant row 3.
$ws.Range("C2").Value = "abc"

# "123" must stay TEXT (not get auto-converted to a number). Build it via
# a formula then paste back as a value so it lands as a shared string
# with no extra number-format/style applied to the cell.
$ws.Range("D2").Formula = '=TEXT(123,"0")'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("E2").Value = "zcx"

$ws.Rows("3:3").Delete()

# Restore the selection state (A2 active, extended over A2:B2).
$ws.Range("A2:B2").Select()
